# "Update dictionary and other changes"
# The "nemad" (stock symbol) column (L) on the Balance Sheet was wrongly
# populated with the company name ("نفت سپاهان"); fix it to show the real
# trading symbol ("شسپا") for every data row (L2:L45). This both adds a
# new shared string and repoints every L-column cell at it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers live in row 1 (L1 = "nemad" header) - leave that alone.
# Rows 2-45 hold the per-period data that should read the ticker symbol.
$ws.Range("L2:L45").Value2 = "شسپا"

# Column L was resized/best-fit to the new (shorter) symbol text.
# ColumnWidth is expressed in characters; the file's internal column width
# unit runs 5/6 of a character wider, so 9.1666... (=55/6) on the COM
# surface serialises to width="10" in the saved workbook.
$ws.Columns("L:L").ColumnWidth = 9.166666666666666

# Leave the cursor parked on L17, matching the saved selection state.
$ws.Range("L17").Select() | Out-Null
